$wb = $excel.ActiveWorkbook

# --- PUTWALL PICKING ---
$ws = $wb.Worksheets.Item("PUTWALL PICKING")
$ws.Range("A2:C200").ClearContents()
$ws.Cells.Item(2,1).Value = "ADOL798N.SEEMANNVAZQ"
$ws.Cells.Item(2,2).Value = 45
$ws.Cells.Item(2,3).Value = 13.98963730569948
$ws.Cells.Item(3,1).Value = "ARJUNBHAI.PATEL"
$ws.Cells.Item(3,2).Value = 27
$ws.Cells.Item(3,3).Value = 8.393782383419689
$ws.Cells.Item(4,1).Value = "BOHD0676.KUSHLIAK"
$ws.Cells.Item(4,2).Value = 161
$ws.Cells.Item(4,3).Value = 50.05181347150259
$ws.Cells.Item(5,1).Value = "DIAN4065.ENTRIALGO"
$ws.Cells.Item(5,2).Value = 27
$ws.Cells.Item(5,3).Value = 8.393782383419689
$ws.Cells.Item(6,1).Value = "LOANA.MBONGO"
$ws.Cells.Item(6,2).Value = 48
$ws.Cells.Item(6,3).Value = 14.92227979274611
$ws.Cells.Item(7,1).Value = "MICA0432.RIZKALLAMAR"
$ws.Cells.Item(7,2).Value = 52
$ws.Cells.Item(7,3).Value = 16.16580310880829
$ws.Cells.Item(8,1).Value = "PATR5027.AMEH"
$ws.Cells.Item(8,2).Value = 2
$ws.Cells.Item(8,3).Value = 0.6217616580310881
$ws.Cells.Item(9,1).Value = "SURESH.DHAWAN"
$ws.Cells.Item(9,2).Value = 99
$ws.Cells.Item(9,3).Value = 30.77720207253886
$ws.Cells.Item(10,1).Value = "THIE6554.DIALLO"
$ws.Cells.Item(10,2).Value = 106
$ws.Cells.Item(10,3).Value = 32.95336787564766
$ws.Cells.Item(11,1).Value = "TUSHAR.BHATIA"
$ws.Cells.Item(11,2).Value = 23
$ws.Cells.Item(11,3).Value = 7.150259067357513
$ws.Cells.Item(12,1).Value = "ZAHIDGUL.MINHAS"
$ws.Cells.Item(12,2).Value = 5
$ws.Cells.Item(12,3).Value = 1.55440414507772

# --- REGULAR PICK ---
$ws = $wb.Worksheets.Item("REGULAR PICK")
$ws.Range("A2:C200").ClearContents()
$ws.Cells.Item(2,1).Value = "ARJUNBHAI.PATEL"
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = 0.310880829015544
$ws.Cells.Item(3,1).Value = "BOHD0676.KUSHLIAK"
$ws.Cells.Item(3,2).Value = 9
$ws.Cells.Item(3,3).Value = 2.797927461139896
$ws.Cells.Item(4,1).Value = "DIAN4065.ENTRIALGO"
$ws.Cells.Item(4,2).Value = 54
$ws.Cells.Item(4,3).Value = 16.78756476683938
$ws.Cells.Item(5,1).Value = "MARI882N.ABDELKADER"
$ws.Cells.Item(5,2).Value = 8
$ws.Cells.Item(5,3).Value = 2.487046632124352
$ws.Cells.Item(6,1).Value = "PATR5027.AMEH"
$ws.Cells.Item(6,2).Value = 4
$ws.Cells.Item(6,3).Value = 1.243523316062176
$ws.Cells.Item(7,1).Value = "WESL5337.CADETTE"
$ws.Cells.Item(7,2).Value = 30
$ws.Cells.Item(7,3).Value = 9.32642487046632
$ws.Cells.Item(8,1).Value = "ZAHIDGUL.MINHAS"
$ws.Cells.Item(8,2).Value = 4
$ws.Cells.Item(8,3).Value = 1.243523316062176

# --- SINGLE PICK ---
$ws = $wb.Worksheets.Item("SINGLE PICK")
$ws.Range("A2:C200").ClearContents()
$ws.Cells.Item(2,1).Value = "GIGNESH.PATEL"
$ws.Cells.Item(2,2).Value = 35
$ws.Cells.Item(2,3).Value = 10.88082901554404
$ws.Cells.Item(3,1).Value = "KADE3054.ZONGO"
$ws.Cells.Item(3,2).Value = 14
$ws.Cells.Item(3,3).Value = 4.352331606217616
$ws.Cells.Item(4,1).Value = "LOANA.MBONGO"
$ws.Cells.Item(4,2).Value = 50
$ws.Cells.Item(4,3).Value = 15.5440414507772
$ws.Cells.Item(5,1).Value = "SEPIDEH.AZARIHASHJIN"
$ws.Cells.Item(5,2).Value = 153
$ws.Cells.Item(5,3).Value = 47.56476683937824
$ws.Cells.Item(6,1).Value = "STAN9294.BAUER"
$ws.Cells.Item(6,2).Value = 86
$ws.Cells.Item(6,3).Value = 26.73575129533679
$ws.Cells.Item(7,1).Value = "TUSHAR.BHATIA"
$ws.Cells.Item(7,2).Value = 12
$ws.Cells.Item(7,3).Value = 3.730569948186528
$ws.Cells.Item(8,1).Value = "WESL5337.CADETTE"
$ws.Cells.Item(8,2).Value = 15
$ws.Cells.Item(8,3).Value = 4.66321243523316

# --- REPLENISHMENT PICK ---
$ws = $wb.Worksheets.Item("REPLENISHMENT PICK")
$ws.Range("A2:C200").ClearContents()
$ws.Cells.Item(2,1).Value = "AGNE8120.CARUTH"
$ws.Cells.Item(2,2).Value = 170
$ws.Cells.Item(2,3).Value = 52.84974093264248
$ws.Cells.Item(3,1).Value = "ARJUNBHAI.PATEL"
$ws.Cells.Item(3,2).Value = 82
$ws.Cells.Item(3,3).Value = 25.49222797927461
$ws.Cells.Item(4,1).Value = "BOHD0676.KUSHLIAK"
$ws.Cells.Item(4,2).Value = 23
$ws.Cells.Item(4,3).Value = 7.150259067357513
$ws.Cells.Item(5,1).Value = "BUDD0680.TENNAKOON"
$ws.Cells.Item(5,2).Value = 112
$ws.Cells.Item(5,3).Value = 34.81865284974093
$ws.Cells.Item(6,1).Value = "DEVI789.SINGH"
$ws.Cells.Item(6,2).Value = 87
$ws.Cells.Item(6,3).Value = 27.04663212435233
$ws.Cells.Item(7,1).Value = "DIAN4065.ENTRIALGO"
$ws.Cells.Item(7,2).Value = 91
$ws.Cells.Item(7,3).Value = 28.29015544041451
$ws.Cells.Item(8,1).Value = "GIGNESH.PATEL"
$ws.Cells.Item(8,2).Value = 73
$ws.Cells.Item(8,3).Value = 22.69430051813471
$ws.Cells.Item(9,1).Value = "INUK4091.QAVAVAU"
$ws.Cells.Item(9,2).Value = 96
$ws.Cells.Item(9,3).Value = 29.84455958549223
$ws.Cells.Item(10,1).Value = "JEEW9554.SITUMUDALIG"
$ws.Cells.Item(10,2).Value = 84
$ws.Cells.Item(10,3).Value = 26.1139896373057
$ws.Cells.Item(11,1).Value = "KADE3054.ZONGO"
$ws.Cells.Item(11,2).Value = 50
$ws.Cells.Item(11,3).Value = 15.5440414507772
$ws.Cells.Item(12,1).Value = "LOANA.MBONGO"
$ws.Cells.Item(12,2).Value = 57
$ws.Cells.Item(12,3).Value = 17.72020725388601
$ws.Cells.Item(13,1).Value = "MARI882N.ABDELKADER"
$ws.Cells.Item(13,2).Value = 37
$ws.Cells.Item(13,3).Value = 11.50259067357513
$ws.Cells.Item(14,1).Value = "MICA0432.RIZKALLAMAR"
$ws.Cells.Item(14,2).Value = 103
$ws.Cells.Item(14,3).Value = 32.02072538860104
$ws.Cells.Item(15,1).Value = "NESR2403.ATTALAH"
$ws.Cells.Item(15,2).Value = 76
$ws.Cells.Item(15,3).Value = 23.62694300518135
$ws.Cells.Item(16,1).Value = "OMAR6689.KHAN"
$ws.Cells.Item(16,2).Value = 61
$ws.Cells.Item(16,3).Value = 18.96373056994819
$ws.Cells.Item(17,1).Value = "PRINCE.FORSON"
$ws.Cells.Item(17,2).Value = 101
$ws.Cells.Item(17,3).Value = 31.39896373056995
$ws.Cells.Item(18,1).Value = "STAN9294.BAUER"
$ws.Cells.Item(18,2).Value = 33
$ws.Cells.Item(18,3).Value = 10.25906735751295
$ws.Cells.Item(19,1).Value = "THIE6554.DIALLO"
$ws.Cells.Item(19,2).Value = 85
$ws.Cells.Item(19,3).Value = 26.42487046632124
$ws.Cells.Item(20,1).Value = "TUSHAR.BHATIA"
$ws.Cells.Item(20,2).Value = 99
$ws.Cells.Item(20,3).Value = 30.77720207253886
$ws.Cells.Item(21,1).Value = "WESL5337.CADETTE"
$ws.Cells.Item(21,2).Value = 83
$ws.Cells.Item(21,3).Value = 25.80310880829015
$ws.Cells.Item(22,1).Value = "WILDINE.JEUNE"
$ws.Cells.Item(22,2).Value = 172
$ws.Cells.Item(22,3).Value = 53.47150259067357
$ws.Cells.Item(23,1).Value = "YATI0689.YATIN"
$ws.Cells.Item(23,2).Value = 112
$ws.Cells.Item(23,3).Value = 34.81865284974093
$ws.Cells.Item(24,1).Value = "ZAKI0190.PHILLIPHORS"
$ws.Cells.Item(24,2).Value = 118
$ws.Cells.Item(24,3).Value = 36.68393782383419

# --- QUICK MOVE ---
$ws = $wb.Worksheets.Item("QUICK MOVE")
$ws.Range("A2:C200").ClearContents()
$ws.Cells.Item(2,1).Value = "ADOL798N.SEEMANNVAZQ"
$ws.Cells.Item(2,2).Value = 173
$ws.Cells.Item(2,3).Value = 53.78238341968912
$ws.Cells.Item(3,1).Value = "BOHD0676.KUSHLIAK"
$ws.Cells.Item(3,2).Value = 79
$ws.Cells.Item(3,3).Value = 24.55958549222798
$ws.Cells.Item(4,1).Value = "DEVI789.SINGH"
$ws.Cells.Item(4,2).Value = 33
$ws.Cells.Item(4,3).Value = 10.25906735751295
$ws.Cells.Item(5,1).Value = "DIAN4065.ENTRIALGO"
$ws.Cells.Item(5,2).Value = 154
$ws.Cells.Item(5,3).Value = 47.87564766839378
$ws.Cells.Item(6,1).Value = "ESSE0616.UDEH"
$ws.Cells.Item(6,2).Value = 321
$ws.Cells.Item(6,3).Value = 99.79274611398964
$ws.Cells.Item(7,1).Value = "JEEW9554.SITUMUDALIG"
$ws.Cells.Item(7,2).Value = 10
$ws.Cells.Item(7,3).Value = 3.10880829015544
$ws.Cells.Item(8,1).Value = "MICA0432.RIZKALLAMAR"
$ws.Cells.Item(8,2).Value = 171
$ws.Cells.Item(8,3).Value = 53.16062176165803
$ws.Cells.Item(9,1).Value = "NESR2403.ATTALAH"
$ws.Cells.Item(9,2).Value = 97
$ws.Cells.Item(9,3).Value = 30.15544041450777
$ws.Cells.Item(10,1).Value = "STAN9294.BAUER"
$ws.Cells.Item(10,2).Value = 112
$ws.Cells.Item(10,3).Value = 34.81865284974093
$ws.Cells.Item(11,1).Value = "SURESH.DHAWAN"
$ws.Cells.Item(11,2).Value = 199
$ws.Cells.Item(11,3).Value = 61.86528497409326
$ws.Cells.Item(12,1).Value = "THIE6554.DIALLO"
$ws.Cells.Item(12,2).Value = 86
$ws.Cells.Item(12,3).Value = 26.73575129533679
$ws.Cells.Item(13,1).Value = "WESL5337.CADETTE"
$ws.Cells.Item(13,2).Value = 107
$ws.Cells.Item(13,3).Value = 33.26424870466321
$ws.Cells.Item(14,1).Value = "YATI0689.YATIN"
$ws.Cells.Item(14,2).Value = 42
$ws.Cells.Item(14,3).Value = 13.05699481865285

# --- IDLE TIME ---
$ws = $wb.Worksheets.Item("IDLE TIME")
$ws.Range("A2:B200").ClearContents()
$ws.Cells.Item(2,1).Value = "ADOL798N.SEEMANNVAZQ"
$ws.Cells.Item(2,2).Value = 119
$ws.Cells.Item(3,1).Value = "AGNE8120.CARUTH"
$ws.Cells.Item(3,2).Value = 43
$ws.Cells.Item(4,1).Value = "ARJUNBHAI.PATEL"
$ws.Cells.Item(4,2).Value = 52
$ws.Cells.Item(5,1).Value = "BOHD0676.KUSHLIAK"
$ws.Cells.Item(5,2).Value = 50
$ws.Cells.Item(6,1).Value = "BUDD0680.TENNAKOON"
$ws.Cells.Item(6,2).Value = 47
$ws.Cells.Item(7,1).Value = "DEVI789.SINGH"
$ws.Cells.Item(7,2).Value = 52
$ws.Cells.Item(8,1).Value = "DIAN4065.ENTRIALGO"
$ws.Cells.Item(8,2).Value = 55
$ws.Cells.Item(9,1).Value = "ESSE0616.UDEH"
$ws.Cells.Item(9,2).Value = 109
$ws.Cells.Item(10,1).Value = "GIGNESH.PATEL"
$ws.Cells.Item(10,2).Value = 82
$ws.Cells.Item(11,1).Value = "INUK4091.QAVAVAU"
$ws.Cells.Item(11,2).Value = 120
$ws.Cells.Item(12,1).Value = "JEEW9554.SITUMUDALIG"
$ws.Cells.Item(12,2).Value = 28
$ws.Cells.Item(13,1).Value = "KADE3054.ZONGO"
$ws.Cells.Item(13,2).Value = 24
$ws.Cells.Item(14,1).Value = "LOANA.MBONGO"
$ws.Cells.Item(14,2).Value = 33
$ws.Cells.Item(15,1).Value = "MARI882N.ABDELKADER"
$ws.Cells.Item(15,2).Value = 131
$ws.Cells.Item(16,1).Value = "MICA0432.RIZKALLAMAR"
$ws.Cells.Item(16,2).Value = 26
$ws.Cells.Item(17,1).Value = "NESR2403.ATTALAH"
$ws.Cells.Item(17,2).Value = 60
$ws.Cells.Item(18,1).Value = "OMAR6689.KHAN"
$ws.Cells.Item(18,2).Value = 77
$ws.Cells.Item(19,1).Value = "PATR5027.AMEH"
$ws.Cells.Item(19,2).Value = 162
$ws.Cells.Item(20,1).Value = "PRINCE.FORSON"
$ws.Cells.Item(20,2).Value = 63
$ws.Cells.Item(21,1).Value = "SEPIDEH.AZARIHASHJIN"
$ws.Cells.Item(21,2).Value = 51
$ws.Cells.Item(22,1).Value = "STAN9294.BAUER"
$ws.Cells.Item(22,2).Value = 54
$ws.Cells.Item(23,1).Value = "SURESH.DHAWAN"
$ws.Cells.Item(23,2).Value = 133
$ws.Cells.Item(24,1).Value = "THIE6554.DIALLO"
$ws.Cells.Item(24,2).Value = 46
$ws.Cells.Item(25,1).Value = "TUSHAR.BHATIA"
$ws.Cells.Item(25,2).Value = 68
$ws.Cells.Item(26,1).Value = "WESL5337.CADETTE"
$ws.Cells.Item(26,2).Value = 47
$ws.Cells.Item(27,1).Value = "WILDINE.JEUNE"
$ws.Cells.Item(27,2).Value = 43
$ws.Cells.Item(28,1).Value = "YATI0689.YATIN"
$ws.Cells.Item(28,2).Value = 51
$ws.Cells.Item(29,1).Value = "ZAHIDGUL.MINHAS"
$ws.Cells.Item(29,2).Value = 140
$ws.Cells.Item(30,1).Value = "ZAKI0190.PHILLIPHORS"
$ws.Cells.Item(30,2).Value = 83
# --- Total Units picked by hour ---
$ws = $wb.Worksheets.Item("Total Units picked by hour")
$ws.Cells.Item(2,2).Value = -27
$ws.Cells.Item(2,3).Value = -24
$ws.Cells.Item(2,4).Value = -468
$ws.Cells.Item(2,5).Value = -1
$ws.Cells.Item(3,2).Value = -41
$ws.Cells.Item(3,3).Value = -130
$ws.Cells.Item(3,4).Value = -788
$ws.Cells.Item(3,5).Value = -81
$ws.Cells.Item(4,2).Value = -32
$ws.Cells.Item(4,3).Value = -168
$ws.Cells.Item(4,4).Value = -645
$ws.Cells.Item(4,5).Value = -357
$ws.Cells.Item(5,2).Value = -10
$ws.Cells.Item(5,3).Value = -43
$ws.Cells.Item(5,4).Value = -104
$ws.Cells.Item(5,5).Value = -156
$ws.Cells.Item(6,2).Value = -110
$ws.Cells.Item(6,3).Value = -365
$ws.Cells.Item(6,4).Value = -2005
$ws.Cells.Item(6,5).Value = -595